$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: Year labels 2010 -> 2012
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "2012"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------------
# Sheet1: updated "Average Values" (E column) and "Riddership Effect" (H column)
# data for each factor row (8-15); rows 16-18 keep the same underlying values.
# ---------------------------------------------------------------------------
$ws1.Range("E8").Value = 395946
$ws1.Range("H8").Value = 91796.752888

$ws1.Range("E9").Value = 0
$ws1.Range("H9").Value = 196880.1356

$ws1.Range("E10").Value = 128693.08
$ws1.Range("H10").Value = 33681.5660209

$ws1.Range("E11").Value = 24.04720583
$ws1.Range("H11").Value = 2157.8121798

$ws1.Range("E12").Value = 4.3491
$ws1.Range("H12").Value = -36301.03872099999

$ws1.Range("E13").Value = 22557.55
$ws1.Range("H13").Value = -39230.787926

$ws1.Range("E14").Value = 7.68
$ws1.Range("H14").Value = -10022.390457

$ws1.Range("E15").Value = 7.6
$ws1.Range("H15").Value = -10570.4129851

# New Reporters row: Riddership Effect now explicitly 0 instead of blank
$ws1.Range("H19").Value = 0

# Total Modeled / Observed Ridership totals
$ws1.Range("E20").Value = 1155567.538
$ws1.Range("E21").Value = 1172297

# ---------------------------------------------------------------------------
# Sheet1: formulas in columns G (% Diff) and I (% Diff / Riddership effect
# share) drop the "*100" factor because the cells are now formatted as a
# percentage instead of a raw number.
# ---------------------------------------------------------------------------
for ($r = 8; $r -le 21; $r++) {
    $ws1.Range("G$r").Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
for ($r = 8; $r -le 19; $r++) {
    $ws1.Range("I$r").Formula = "=IFERROR(H$r/`$E`$21,0)"
}

# ---------------------------------------------------------------------------
# Sheet1: number formats -- "Average Values" (E:F), "Riddership Effect" (H)
# become 2-decimal numbers; "% Diff" columns (G, I) become percentages.
# ---------------------------------------------------------------------------
$ws1.Range("E8:F18").NumberFormat = "#,##0.00"
$ws1.Range("H8:H18").NumberFormat = "#,##0.00"
$ws1.Range("G8:G18").NumberFormat = "0.00%"
$ws1.Range("I8:I18").NumberFormat = "0.00%"

$ws1.Range("E19:F19").NumberFormat = "#,##0.00"
$ws1.Range("H19").NumberFormat = "#,##0.00"
$ws1.Range("G19").NumberFormat = "0.00%"
$ws1.Range("I19").NumberFormat = "0.00%"

$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat = "#,##0.00"
$ws1.Range("G20").NumberFormat = "0.00%"
$ws1.Range("I20").NumberFormat = "0.00%"

$ws1.Range("E21:F21").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat = "#,##0.00"
$ws1.Range("G21").NumberFormat = "0.00%"
$ws1.Range("I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Sheet1: refresh the view -- scroll back to the top-left and move the
# selection from K20 to H21.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
$ws1.Range("H21").Select()
